# Peaufinage - Part 4
# Removes spurious spell-check markup (w:proofErr spellStart/spellEnd wrapping
# individual runs) by normalizing the run text back into a single run, and
# fixes a typo "camionnette}" -> "camionnet}".

$d = $word.ActiveDocument

# Merge {adresse}{cp}{ville} back into a single run of text (removes the
# proofErr-wrapped "cp" run split). Occurs twice (main paragraph + table copy).
$d.Content.Find.Execute("{adresse}{cp}{ville}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{adresse}{cp}{ville}", 2)

# Merge {siret} back into a single run of text. Occurs twice.
$d.Content.Find.Execute("{siret}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{siret}", 2)

# Merge {activite} back into a single run of text. Occurs twice.
$d.Content.Find.Execute("{activite}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{activite}", 2)

# Fix typo: camionnette} -> camionnet}
$d.Content.Find.Execute("camionnette}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "camionnet}", 2)

# Merge {deuxroues} back into a single run of text.
$d.Content.Find.Execute("{deuxroues}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{deuxroues}", 2)

# Merge "(sez) " back into a single run of text.
$d.Content.Find.Execute("(sez) ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "(sez) ", 2)
